# Case_2_129 res_bus vm_pu.xlsx update
# Updates bus voltage magnitude (vm_pu) results for the 380 kV slack-bus case
# (slack bus voltage set-point changed from 1.05 pu to 1.02 pu, rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vmPuData = @{
    2 = @{ "B"=1.02; "C"=1.016352683689298; "D"=1.022454944189424; "E"=0.9926147277508489; "F"=1.014695902335164; "I"=1.026630963468057; "J"=1.021572980513314; "K"=1.025289275695335; "L"=0.9955398523336033; "M"=1.01755328198371; "N"=1.023023731321702 }
    3 = @{ "B"=1.02; "C"=1.017417861839457; "D"=1.023248169251022; "E"=0.9936372048519304; "F"=1.016405834474572; "I"=1.026827765241115; "J"=1.022272963540431; "K"=1.025889207207821; "L"=0.9963617723202692; "M"=1.019065705187147; "N"=1.023724708405006 }
    4 = @{ "B"=1.02; "C"=1.018105911866639; "D"=1.023759904275757; "E"=0.9942998659930995; "F"=1.017510845966939; "I"=1.026952620334731; "J"=1.022724194184579; "K"=1.026275222182209; "L"=0.9968940712668345; "M"=1.020042484989629; "N"=1.024176579848425 }
    5 = @{ "B"=1.02; "C"=1.018394885311058; "D"=1.023974671223621; "E"=0.9945786998346017; "F"=1.017975060115984; "I"=1.027004513831203; "J"=1.02291348499487; "K"=1.026436981812285; "L"=0.997117960005301; "M"=1.020452686370934; "N"=1.024366139473379 }
    6 = @{ "B"=1.02; "C"=1.018443388711805; "D"=1.024010710076606; "E"=0.9946255319796338; "F"=1.018052984477822; "I"=1.027013192058867; "J"=1.022945243962417; "K"=1.026464111439777; "L"=0.9971555583673453; "M"=1.02052153560906; "N"=1.024397943542303 }
    7 = @{ "B"=1.02; "C"=1.018109774252115; "D"=1.02376277544106; "E"=0.9943035907982488; "F"=1.017517050114743; "I"=1.026953316078162; "J"=1.022726725093109; "K"=1.026277385668527; "L"=0.9968970624462087; "M"=1.020047967826948; "N"=1.024179114351136 }
    8 = @{ "B"=1.02; "C"=1.016712913817593; "D"=1.022723336642958; "E"=0.9929600610674301; "F"=1.015274083460961; "I"=1.026697989025844; "J"=1.021809897138386; "K"=1.025492477543771; "L"=0.995817528259106; "M"=1.018064802407768; "N"=1.023260984395556 }
    9 = @{ "B"=1.02; "C"=1.014242230554431; "D"=1.020879915135877; "E"=0.9906006454969559; "F"=1.011310322327099; "I"=1.026228996634949; "J"=1.020181193002234; "K"=1.024092618545396; "L"=0.9939188001724441; "M"=1.014555602028238; "N"=1.021629967312716 }
    10 = @{ "B"=1.02; "C"=1.012588738015894; "D"=1.019642968874427; "E"=0.989033133672735; "F"=1.008659577465204; "I"=1.025903491092794; "J"=1.019086447570949; "K"=1.023148045563451; "L"=0.9926553831429383; "M"=1.012205791179569; "N"=1.020533667217349 }
    11 = @{ "B"=1.02; "C"=1.011871209637478; "D"=1.019105444021029; "E"=0.988355674866747; "F"=1.007509688740975; "I"=1.025759491817609; "J"=1.018610263879704; "K"=1.022736329627008; "L"=0.9921088820399291; "M"=1.011185726843843; "N"=1.020056807290644 }
    12 = @{ "B"=1.02; "C"=1.01160445068494; "D"=1.018905493495206; "E"=0.9881042295826724; "F"=1.007082242005985; "I"=1.025705544789199; "J"=1.018433062444378; "K"=1.022582991282056; "L"=0.9919059725120875; "M"=1.010806431639418; "N"=1.019879354208954 }
    13 = @{ "B"=1.02; "C"=1.011661682162603; "D"=1.01894839670407; "E"=0.9881581567098651; "F"=1.00717394575956; "I"=1.025717137403012; "J"=1.018471087526327; "K"=1.022615901399273; "L"=0.9919494934313052; "M"=1.010887809917943; "N"=1.019917433290881 }
    14 = @{ "B"=1.02; "C"=1.011849164112576; "D"=1.019088921970454; "E"=0.9883348863814464; "F"=1.007474362595664; "I"=1.025755041912102; "J"=1.018595623009816; "K"=1.022723662978779; "L"=0.9920921077337197; "M"=1.011154382359933; "N"=1.020042145629041 }
    15 = @{ "B"=1.02; "C"=1.011964646429646; "D"=1.019175465768117; "E"=0.9884438009545853; "F"=1.007659415663257; "I"=1.025778335253778; "J"=1.018672310214152; "K"=1.022790004236537; "L"=0.9921799884222134; "M"=1.011318573409226; "N"=1.020118941738003 }
    16 = @{ "B"=1.02; "C"=1.012636324963105; "D"=1.019678602066127; "E"=0.9890781214508737; "F"=1.008735846591285; "I"=1.025912983474025; "J"=1.019118004781412; "K"=1.023175312567926; "L"=0.9926916645766087; "M"=1.01227343425134; "N"=1.020565269242671 }
    17 = @{ "B"=1.02; "C"=1.01305723246406; "D"=1.01999369132389; "E"=0.9894763578477731; "F"=1.009410493262613; "I"=1.025996626911626; "J"=1.01939699949588; "K"=1.023416279729347; "L"=0.9930127773692701; "M"=1.012871695146693; "N"=1.020844660161634 }
    18 = @{ "B"=1.02; "C"=1.013302590861348; "D"=1.020177292514019; "E"=0.9897087662937551; "F"=1.009803801543454; "I"=1.026045120162039; "J"=1.01955952494631; "K"=1.023556570473386; "L"=0.9932001317071766; "M"=1.013220402293914; "N"=1.021007416416846 }
    19 = @{ "B"=1.02; "C"=1.013386226420138; "D"=1.02023986443479; "E"=0.9897880325774039; "F"=1.009937875637839; "I"=1.026061605162946; "J"=1.019614906808336; "K"=1.023604361714805; "L"=0.993264023964098; "M"=1.013339260607731; "N"=1.021062876927468 }
    20 = @{ "B"=1.02; "C"=1.013012088595076; "D"=1.019959904380522; "E"=0.9894336180360677; "F"=1.009338131009134; "I"=1.025987683233767; "J"=1.01936708747835; "K"=1.023390453275968; "L"=0.9929783193494215; "M"=1.012807533174147; "N"=1.020814705665609 }
    21 = @{ "B"=1.02; "C"=1.011793961919384; "D"=1.019047548796336; "E"=0.9882828385668249; "F"=1.007385906441489; "I"=1.025743892662174; "J"=1.018558959421325; "K"=1.022691941175273; "L"=0.9920501090198102; "M"=1.011075894485766; "N"=1.020005429974049 }
    22 = @{ "B"=1.02; "C"=1.011026704654397; "D"=1.018472236776189; "E"=0.9875604150241495; "F"=1.006156568441006; "I"=1.025587954266781; "J"=1.018048972170091; "K"=1.022250393365223; "L"=0.9914670000341481; "M"=1.009984837334178; "N"=1.019494718482421 }
    23 = @{ "B"=1.02; "C"=1.011433573479443; "D"=1.018777380136001; "E"=0.9879432794643023; "F"=1.006808447469417; "I"=1.025670872266445; "J"=1.018319505535985; "K"=1.022484690991812; "L"=0.991776070289318; "M"=1.0105634494877; "N"=1.019765636036726 }
    24 = @{ "B"=1.02; "C"=1.013032487629076; "D"=1.019975171817649; "E"=0.9894529299347244; "F"=1.009370829018875; "I"=1.02599172540694; "J"=1.019380604072392; "K"=1.023402123945975; "L"=0.9929938892766442; "M"=1.01283652597404; "N"=1.020828241454765 }
    25 = @{ "B"=1.02; "C"=1.014882072341768; "D"=1.021357888565543; "E"=0.9912096547607049; "F"=1.012336457598084; "I"=1.026352504695465; "J"=1.020603820029821; "K"=1.024456507912286; "L"=0.9939188001724441; "M"=1.014555602028238; "N"=1.021629967312716 }
}

foreach ($rowNum in $vmPuData.Keys) {
    $rowData = $vmPuData[$rowNum]
    foreach ($colLetter in $rowData.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $rowData[$colLetter]
    }
}
